$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I113").Value = 3800.8  # was 3852.5
$ws.Range("H113").Value = 6834  # was 9901.666999999999
$ws.Range("M113").Value = -546.8000000000002  # was -598.5
$ws.Range("K113").Value = 3800.8  # was 3852.5
$ws.Range("L131").Value = 11896.6671  # was 12145.5
$ws.Range("H131").Value = 3392.647  # was 3470.5557
$ws.Range("N131").Value = -21976.6671  # was -22225.5
$ws.Range("J131").Value = 3965.5557  # was 4048.5
$ws.Range("I132").Value = 29415828  # was 32262564
$ws.Range("M132").Value = -88244954  # was -96785162
$ws.Range("L132").Value = 19800  # was 22500
$ws.Range("K132").Value = 88247484  # was 96787692
$ws.Range("N132").Value = -24860  # was -27560
$ws.Range("J132").Value = 6600  # was 7500
$ws.Range("H132").Value = 25645414  # was 28576270
$ws.Range("J138").Value = 2861.244  # was 2851.6707
$ws.Range("K138").Value = 1899.3333  # was 1958.1177
$ws.Range("N138").Value = -18863.732  # was -18835.0121
$ws.Range("M138").Value = 3240.6667  # was 3181.8823
$ws.Range("L138").Value = 8583.732  # was 8555.0121
$ws.Range("H138").Value = 2460.18  # was 2474.0708
$ws.Range("I138").Value = 633.1111  # was 652.7059
$ws.Range("L140").Value = 46266.875  # was 46675.805
$ws.Range("J140").Value = 46266.875  # was 46675.805
$ws.Range("N140").Value = -56626.875  # was -57035.805
$ws.Range("H140").Value = 45692.35  # was 46059.09

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M4").Value = -917.3334  # was -834.75
$ws.Range("I4").Value = 1033.3334  # was 950.75
$ws.Range("K4").Value = 1033.3334  # was 950.75
$ws.Range("H4").Value = 1033.3334  # was 950.75
$ws.Range("J6").Value = 14798  # was 14798.333
$ws.Range("I6").Value = 3200  # was 500
$ws.Range("N6").Value = -15144  # was -15144.333
$ws.Range("M6").Value = -3027  # was -327
$ws.Range("H6").Value = 10932  # was 11223.75
$ws.Range("K6").Value = 3200  # was 500
$ws.Range("L6").Value = 14798  # was 14798.333
$ws.Range("K22").Value = 2500  # was 4000
$ws.Range("L22").Value = 0  # was 9750
$ws.Range("I22").Value = 2500  # was 4000
$ws.Range("N22").ClearContents()  # was -10348
$ws.Range("M22").Value = -2201  # was -3701
$ws.Range("H22").Value = 2500  # was 7833.3335
$ws.Range("J22").Value = 0  # was 9750
$ws.Range("H139").Value = 42579  # was 43889.523
$ws.Range("N139").Value = -52859  # was -54169.523
$ws.Range("J139").Value = 42579  # was 43889.523
$ws.Range("L139").Value = 42579  # was 43889.523

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K24").Value = 3250  # was 4000
$ws.Range("H24").Value = 3250  # was 4000
$ws.Range("M24").Value = -3015  # was -3765
$ws.Range("I24").Value = 3250  # was 4000

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L31").Value = 7501.625  # was 6745.8887
$ws.Range("N31").Value = -8091.625  # was -7335.8887
$ws.Range("J31").Value = 7501.625  # was 6745.8887
$ws.Range("M31").Value = -760  # was -743.7221999999999
$ws.Range("K31").Value = 1055  # was 1038.7222
$ws.Range("H31").Value = 3117.92  # was 2941.111
$ws.Range("I31").Value = 1055  # was 1038.7222
$ws.Range("J34").Value = 7501.625  # was 6745.8887
$ws.Range("I34").Value = 1055  # was 1038.7222
$ws.Range("L34").Value = 7501.625  # was 6745.8887
$ws.Range("K34").Value = 1055  # was 1038.7222
$ws.Range("H34").Value = 3117.92  # was 2941.111
$ws.Range("N34").Value = -7905.625  # was -7149.8887
$ws.Range("M34").Value = -853  # was -836.7221999999999
$ws.Range("I58").Value = 1854.6957  # was 1831.1915
$ws.Range("M58").Value = -1651.6957  # was -1628.1915
$ws.Range("H58").Value = 3219.9285  # was 3176.5964
$ws.Range("K58").Value = 1854.6957  # was 1831.1915
$ws.Range("L99").Value = 6100  # was 6580
$ws.Range("N99").Value = -9096  # was -9576
$ws.Range("H99").Value = 10530505  # was 11115461
$ws.Range("J99").Value = 6100  # was 6580
$ws.Range("N126").Value = -23240  # was -24680
$ws.Range("L126").Value = 18300  # was 19740
$ws.Range("H126").Value = 10530505  # was 11115461
$ws.Range("J126").Value = 6100  # was 6580
$ws.Range("I132").Value = 1745.5  # was 1736.1538
$ws.Range("M132").Value = -2706.5  # was -2678.4614
$ws.Range("L132").Value = 18598.8  # was 22498.9995
$ws.Range("K132").Value = 5236.5  # was 5208.4614
$ws.Range("N132").Value = -23658.8  # was -27558.9995
$ws.Range("J132").Value = 6199.6  # was 7499.6665
$ws.Range("H132").Value = 3055.5293  # was 2816.8125
$ws.Range("K136").Value = 5564.0871  # was 5493.5745
$ws.Range("M136").Value = -3014.0871  # was -2943.5745
$ws.Range("I136").Value = 1854.6957  # was 1831.1915
$ws.Range("H136").Value = 3219.9285  # was 3176.5964

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L22").Value = 5502  # was 5469.3333
$ws.Range("N22").Value = -5840  # was -5807.3333
$ws.Range("H22").Value = 1834  # was 1823.1111
$ws.Range("J22").Value = 1834  # was 1823.1111
$ws.Range("L27").Value = 5502  # was 5469.3333
$ws.Range("H27").Value = 1834  # was 1823.1111
$ws.Range("N27").Value = -5706  # was -5673.3333
$ws.Range("J27").Value = 1834  # was 1823.1111
$ws.Range("H68").Value = 26099.5  # was 26049.75
$ws.Range("K68").Value = 3597  # was 3000
$ws.Range("I68").Value = 1199  # was 1000
$ws.Range("M68").Value = -2786  # was -2189
$ws.Range("H71").Value = 26099.5  # was 26049.75
$ws.Range("K71").Value = 10791  # was 9000
$ws.Range("M71").Value = -6735  # was -4944
$ws.Range("I71").Value = 1199  # was 1000
$ws.Range("N113").Value = -26791745  # was -25006038.5
$ws.Range("I113").Value = 566.5789  # was 590.8823
$ws.Range("H113").Value = 3788444.2  # was 3906829.2
$ws.Range("M113").Value = 470.2633000000001  # was 397.3531
$ws.Range("L113").Value = 26787405  # was 25001698.5
$ws.Range("J113").Value = 8929135  # was 8333899.5
$ws.Range("K113").Value = 1699.7367  # was 1772.6469
$ws.Range("M139").Value = 2543.5  # was 1810
$ws.Range("H139").Value = 1196.0869  # was 1762.3182
$ws.Range("N139").Value = -20480  # was -18991.625
$ws.Range("J139").Value = 3400  # was 2903.875
$ws.Range("K139").Value = 2596.5  # was 3330
$ws.Range("I139").Value = 865.5  # was 1110
$ws.Range("L139").Value = 10200  # was 8711.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J24").Value = 0  # was 9000
$ws.Range("L24").Value = 0  # was 9000
$ws.Range("H24").Value = 0  # was 9000
$ws.Range("N24").ClearContents()  # was -9346
$ws.Range("I97").Value = 1011.61536  # was 960.73334
$ws.Range("M97").Value = -515.61536  # was -464.73334
$ws.Range("K97").Value = 1011.61536  # was 960.73334
$ws.Range("H97").Value = 1008.8125  # was 966.7222
$ws.Range("I132").Value = 2142.8572  # was 2400
$ws.Range("M132").Value = -3898.571599999999  # was -4670
$ws.Range("L132").Value = 20485.2861  # was 21249.4995
$ws.Range("K132").Value = 6428.571599999999  # was 7200
$ws.Range("N132").Value = -25545.2861  # was -26309.4995
$ws.Range("J132").Value = 6828.4287  # was 7083.1665
$ws.Range("H132").Value = 4485.643  # was 4954.4546

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J6").Value = 40000  # was 0
$ws.Range("N6").Value = -40224  # was empty
$ws.Range("H6").Value = 40000  # was 0
$ws.Range("L6").Value = 40000  # was 0
$ws.Range("J128").Value = 41900  # was 41993.332
$ws.Range("L128").Value = 41900  # was 41993.332
$ws.Range("H128").Value = 41900  # was 41993.332
$ws.Range("N128").Value = -51860  # was -51953.332

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L20").Value = 11000  # was 0
$ws.Range("N20").Value = -11480  # was empty
$ws.Range("H20").Value = 11000  # was 0
$ws.Range("J20").Value = 11000  # was 0
$ws.Range("J24").Value = 17495  # was 16247.25
$ws.Range("L24").Value = 17495  # was 16247.25
$ws.Range("H24").Value = 17495  # was 16247.25
$ws.Range("N24").Value = -17955  # was -16707.25
$ws.Range("I30").Value = 16959.25  # was 9500
$ws.Range("N30").Value = -28133.8  # was -11214
$ws.Range("M30").Value = -16852.25  # was -9393
$ws.Range("J30").Value = 27919.8  # was 11000
$ws.Range("L30").Value = 27919.8  # was 11000
$ws.Range("H30").Value = 23048.445  # was 10400
$ws.Range("K30").Value = 16959.25  # was 9500
$ws.Range("L99").Value = 35000  # was 0
$ws.Range("N99").Value = -40990  # was empty
$ws.Range("K99").Value = 25000  # was 0
$ws.Range("I99").Value = 25000  # was 0
$ws.Range("M99").Value = -22005  # was empty
$ws.Range("H99").Value = 30000  # was 0
$ws.Range("J99").Value = 35000  # was 0
$ws.Range("N113").Value = -8094.5  # was -8840
$ws.Range("H113").Value = 1072.6  # was 1090
$ws.Range("L113").Value = 3754.5  # was 4500
$ws.Range("J113").Value = 1251.5  # was 1500
$ws.Range("N126").Value = -7996637.600000001  # was -22239.9995
$ws.Range("I126").Value = 1844  # was 1552.6666
$ws.Range("M126").Value = -3062  # was -2187.9998
$ws.Range("L126").Value = 7991697.600000001  # was 17299.9995
$ws.Range("K126").Value = 5532  # was 4657.9998
$ws.Range("H126").Value = 969864.0600000001  # was 2957.3333
$ws.Range("J126").Value = 2663899.2  # was 5766.6665
